# Insert two new data rows at the top of the Ajo price-history block (rows
# 738-739), pushing the existing rows 738-761 down to 740-763.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("738:739").Insert()

# New row 738
$ws.Cells.Item(738, 1).Value = 3
$ws.Cells.Item(738, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(738, 3).Value = "Coquimbo"
$ws.Cells.Item(738, 4).Value = 45075
$ws.Cells.Item(738, 5).Value = 5
$ws.Cells.Item(738, 6).Value = 100112003
$ws.Cells.Item(738, 7).Value = "Ajo"
$ws.Cells.Item(738, 8).Value = "Chino"
$ws.Cells.Item(738, 9).Value = "Primera"
$ws.Cells.Item(738, 10).Value = 90
$ws.Cells.Item(738, 11).Value = 15000
$ws.Cells.Item(738, 12).Value = 15500
$ws.Cells.Item(738, 13).Value = 15250
$ws.Cells.Item(738, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(738, 15).Value = "China"
$ws.Cells.Item(738, 16).Value = 1525
$ws.Cells.Item(738, 17).Value = 10
$ws.Cells.Item(738, 18).Value = "Hortaliza"

# New row 739
$ws.Cells.Item(739, 1).Value = 3
$ws.Cells.Item(739, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(739, 3).Value = "Coquimbo"
$ws.Cells.Item(739, 4).Value = 45075
$ws.Cells.Item(739, 5).Value = 5
$ws.Cells.Item(739, 6).Value = 100112003
$ws.Cells.Item(739, 7).Value = "Ajo"
$ws.Cells.Item(739, 8).Value = "Chino"
$ws.Cells.Item(739, 9).Value = "Primera"
$ws.Cells.Item(739, 10).Value = 105
$ws.Cells.Item(739, 11).Value = 16500
$ws.Cells.Item(739, 12).Value = 17000
$ws.Cells.Item(739, 13).Value = 16738
$ws.Cells.Item(739, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(739, 15).Value = "China"
$ws.Cells.Item(739, 16).Value = 1674
$ws.Cells.Item(739, 17).Value = 10
$ws.Cells.Item(739, 18).Value = "Hortaliza"

# Match the date format used by the rest of column D (custom date/time
# number format carried by the existing cells, e.g. D740).
$ws.Range("D738:D739").NumberFormat = $ws.Range("D740").NumberFormat()
